# Weekly update: insert a new price record for Vega Modelo de Temuco - Mango
# as row 231, shifting the existing rows 231-241 down to 232-242.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 231 (pushes old rows 231..241 -> 232..242)
$ws.Rows(231).Insert()

# Populate the newly inserted row with the new weekly record
$ws.Range("A231").Value = 10
$ws.Range("B231").Value = "Vega Modelo de Temuco"
$ws.Range("C231").Value = "La Araucanía"
$ws.Range("D231").Value = 44509
$ws.Range("E231").Value = 9
$ws.Range("F231").Value = "Fruta"
$ws.Range("G231").Value = 100108
$ws.Range("H231").Value = "Tropicales y subtropicales"
$ws.Range("I231").Value = 100108002
$ws.Range("J231").Value = "Mango"
$ws.Range("K231").Value = "Sin especificar"
$ws.Range("L231").Value = "Primera"
$ws.Range("M231").Value = 210
$ws.Range("N231").Value = 8000
$ws.Range("O231").Value = 8000
$ws.Range("P231").Value = 8000
$ws.Range("Q231").Value = "$/bandeja 4 kilos"
$ws.Range("R231").Value = "Perú"
$ws.Range("S231").Value = 2000
$ws.Range("T231").Value = 4
